# Apply the "Updated symbol list" refresh: new Price/Volume(1h) figures
# and the Hora column bumped from 11 to 12 across every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='303.03'},
    @{Cell='G2'; Value='12'},
    @{Cell='D3'; Value='32.02'},
    @{Cell='E3'; Value='8.84%'},
    @{Cell='G3'; Value='12'},
    @{Cell='D4'; Value='5.216'},
    @{Cell='E4'; Value='2.92%'},
    @{Cell='G4'; Value='12'},
    @{Cell='D5'; Value='0.07384'},
    @{Cell='E5'; Value='9.70%'},
    @{Cell='G5'; Value='12'},
    @{Cell='D6'; Value='7.830'},
    @{Cell='E6'; Value='6.76%'},
    @{Cell='G6'; Value='12'},
    @{Cell='D7'; Value='3.736'},
    @{Cell='E7'; Value='8.24%'},
    @{Cell='G7'; Value='12'},
    @{Cell='D8'; Value='1.471'},
    @{Cell='E8'; Value='6.57%'},
    @{Cell='G8'; Value='12'},
    @{Cell='D9'; Value='0.9053'},
    @{Cell='E9'; Value='0.10%'},
    @{Cell='G9'; Value='12'},
    @{Cell='D10'; Value='0.01667'},
    @{Cell='E10'; Value='2,477.85%'},
    @{Cell='G10'; Value='12'},
    @{Cell='D11'; Value='0.1682'},
    @{Cell='E11'; Value='5.82%'},
    @{Cell='G11'; Value='12'},
    @{Cell='D12'; Value='0.07449'},
    @{Cell='E12'; Value='7.03%'},
    @{Cell='G12'; Value='12'},
    @{Cell='D13'; Value='0.08009'},
    @{Cell='E13'; Value='5.36%'},
    @{Cell='G13'; Value='12'},
    @{Cell='D14'; Value='0.03041'},
    @{Cell='E14'; Value='3.82%'},
    @{Cell='G14'; Value='12'},
    @{Cell='D15'; Value='0.09900'},
    @{Cell='E15'; Value='10.11%'},
    @{Cell='G15'; Value='12'},
    @{Cell='D16'; Value='0.001512'},
    @{Cell='E16'; Value='-4.90%'},
    @{Cell='G16'; Value='12'},
    @{Cell='D17'; Value='0.04542'},
    @{Cell='E17'; Value='1.38%'},
    @{Cell='G17'; Value='12'},
    @{Cell='D18'; Value='0.006276'},
    @{Cell='E18'; Value='0.47%'},
    @{Cell='G18'; Value='12'},
    @{Cell='D19'; Value='3.487'},
    @{Cell='E19'; Value='1.13%'},
    @{Cell='G19'; Value='12'},
    @{Cell='D20'; Value='2.232'},
    @{Cell='E20'; Value='0.15%'},
    @{Cell='G20'; Value='12'},
    @{Cell='D21'; Value='0.3338'},
    @{Cell='E21'; Value='4.16%'},
    @{Cell='G21'; Value='12'},
    @{Cell='D22'; Value='0.1334'},
    @{Cell='E22'; Value='1.42%'},
    @{Cell='G22'; Value='12'},
    @{Cell='D23'; Value='4.500'},
    @{Cell='E23'; Value='11.18%'},
    @{Cell='G23'; Value='12'},
    @{Cell='D24'; Value='0.1648'},
    @{Cell='E24'; Value='4.26%'},
    @{Cell='G24'; Value='12'},
    @{Cell='D25'; Value='0.001215'},
    @{Cell='E25'; Value='1.70%'},
    @{Cell='G25'; Value='12'},
    @{Cell='D26'; Value='0.004433'},
    @{Cell='E26'; Value='1.37%'},
    @{Cell='G26'; Value='12'},
    @{Cell='D27'; Value='0.0001298'},
    @{Cell='E27'; Value='8.12%'},
    @{Cell='G27'; Value='12'},
    @{Cell='D28'; Value='0.0001740'},
    @{Cell='E28'; Value='7.44%'},
    @{Cell='G28'; Value='12'},
    @{Cell='G29'; Value='12'},
    @{Cell='G30'; Value='12'},
    @{Cell='G31'; Value='12'},
    @{Cell='G32'; Value='12'},
    @{Cell='G33'; Value='12'},
    @{Cell='G34'; Value='12'},
    @{Cell='G35'; Value='12'},
    @{Cell='G36'; Value='12'},
    @{Cell='G37'; Value='12'},
    @{Cell='G38'; Value='12'},
    @{Cell='G39'; Value='12'},
    @{Cell='D40'; Value='0.04486'},
    @{Cell='E40'; Value='5.43%'},
    @{Cell='G40'; Value='12'},
    @{Cell='D41'; Value='0.007212'},
    @{Cell='E41'; Value='5.41%'},
    @{Cell='G41'; Value='12'},
    @{Cell='D42'; Value='0.1346'},
    @{Cell='E42'; Value='8.55%'},
    @{Cell='G42'; Value='12'},
    @{Cell='D43'; Value='0.002386'},
    @{Cell='E43'; Value='6.96%'},
    @{Cell='G43'; Value='12'},
    @{Cell='D44'; Value='0.01404'},
    @{Cell='E44'; Value='10.75%'},
    @{Cell='G44'; Value='12'},
    @{Cell='D45'; Value='0.00006152'},
    @{Cell='E45'; Value='7.84%'},
    @{Cell='G45'; Value='12'},
    @{Cell='D46'; Value='0.7068'},
    @{Cell='E46'; Value='-64.04%'},
    @{Cell='G46'; Value='12'},
    @{Cell='D47'; Value='0.01299'},
    @{Cell='E47'; Value='-13.45%'},
    @{Cell='G47'; Value='12'},
    @{Cell='G48'; Value='12'},
    @{Cell='G49'; Value='12'},
    @{Cell='G50'; Value='12'},
    @{Cell='G51'; Value='12'}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    # Force a text number format before writing so Excel keeps the value
    # as a literal string (matching the original inline-string cells)
    # instead of silently coercing numeric-looking text to a Number.
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    # Reset the style back to Normal so we do not leave a stray
    # "Text format" style attached to the cell (the source workbook
    # cells carry no explicit style index).
    $r.Style = "Normal"
}
